$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Immanuel Quickley", "PG,SG",    "Toronto Raptors"),
    @("Desmond Bane",       "SG,SF",    "Memphis Grizzlies"),
    @("Jalen Brunson",      "PG",       "New York Knicks"),
    @("Norman Powell",      "SG,SF",    "LA Clippers"),
    @("Devin Vassell",      "SG,SF",    "San Antonio Spurs"),
    @("Kawhi Leonard",      "SG,SF,PF", "LA Clippers"),
    @("LeBron James",       "SF,PF",    "Los Angeles Lakers"),
    @("Tim Hardaway Jr.",   "SG,SF",    "Detroit Pistons"),
    @("Myles Turner",       "C",        "Indiana Pacers"),
    @("P.J. Washington",    "PF",       "Dallas Mavericks"),
    @("Jamal Murray",       "PG,SG",    "Denver Nuggets"),
    @("Devin Booker",       "PG,SG",    "Phoenix Suns"),
    @("Ben Simmons",        "PG,C",     "Brooklyn Nets"),
    @("Coby White",         "PG,SG",    "Chicago Bulls"),
    @("Trae Young",         "PG",       "Atlanta Hawks"),
    @("D'Angelo Russell",   "PG",       "Brooklyn Nets"),
    @("Brandon Ingram",     "SG,SF,PF", "New Orleans Pelicans"),
    @("Walker Kessler",     "C",        "Utah Jazz")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
